# update gsi indicator for 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 = Western Gulf Stream Index indicator.
# Update the "Status in 2024" column (B) with the refreshed long-term-average
# comparison, and the "Time Series" column (D) with the refreshed chart image
# filename for the 2025 update.
$ws.Range("B5").Value = "Near long term (1996-2025) average"
$ws.Range("D5").Value = "western gulf stream index_2026-02-23.png"

# Leave the selection on the cell that was last edited, matching the saved
# workbook's cursor position.
$ws.Range("D5").Select()
